# Fix the bug-tracker row about the "previous track" button: expand the
# root-cause note in column C (row 20) with the actual fix description, and
# mark the bug as resolved ("√") in column D, same as the other rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 20 ("实在按上一曲的之后..."): append the root-cause / fix explanation.
$ws.Range("C20").Value = "实在按上一曲的之后会操作导致相反，如果按发送重播指令正常，通过将上一曲的指令改为setsource+playback直接播放而不是replay后问题解决，我觉得可能是player setSource之后本来就没再播放状态，此时调用stop可能会出问题。"

# Mark the bug resolved in column D (same "√" shared string used elsewhere).
$ws.Range("D20").Value = "√"

# The note is much longer now, so grow the row to fit it.
$ws.Rows.Item(20).RowHeight = 82

# Update the window scroll position / selection to where the edit happened.
$ws.Range("F22").Select()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
